# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 654
    $ws.Range("F3").Value = 3895
    $ws.Range("F5").Value = 730
}
